# Fix some disprop indices and fully report them in codebook.
# The Schleswig-Holstein government that was previously marked as "ongoing"
# (G=401768) actually ended on 2022-06-29 (44741), and a new ("ongoing")
# government (Günther II, CDU~~~Grüne) needs to be inserted right after it,
# which shifts the whole Thüringen block (old rows 361-370) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 361 (pushes old rows 361-370 -> 362-371),
# inheriting number formats/styles from the row above like Excel normally does.
$ws.Rows.Item(361).Insert() | Out-Null

# Row 360 (Albig -> Günther handover): the old government's end date was a
# placeholder "still serving" sentinel; now it has a real end date.
$ws.Range("G360").Value = 44741

# New row 361: Schleswig-Holstein, Günther I cabinet (CDU~~~Grüne).
$ws.Range("A361").Value = 11630
$ws.Range("B361").Value = "Schleswig-Holstein"
$ws.Range("C361").Value = 44689
$ws.Range("E361").Value = 30
$ws.Range("F361").Value = 44741
$ws.Range("G361").Value = 401768
$ws.Range("H361").Value = "CDU~~~Grüne"
$ws.Range("I361").Value = "Günther, Daniel"
$ws.Range("J361").Value = "CDU"

# Reflect the new scroll position / active cell from the authored view.
$win = $excel.ActiveWindow
$win.ScrollRow = 325
$win.ScrollColumn = 1
$ws.Range("A362").Select()
